$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 630; this shifts the existing rows 630-675
# down to 631-676, growing the used range from A1:R675 to A1:R676.
$ws.Rows.Item(630).Insert()

# Populate the newly inserted row 630 with the new data record.
$ws.Range("A630").Value = 3
$ws.Range("B630").Value = "Femacal de La Calera"
$ws.Range("C630").Value = "Coquimbo"
$ws.Range("D630").Value = 45021
$ws.Range("E630").Value = 5
$ws.Range("F630").Value = 100112037
$ws.Range("G630").Value = "Cebollín"
$ws.Range("H630").Value = "Sin especificar"
$ws.Range("I630").Value = "Primera"
$ws.Range("J630").Value = 220
$ws.Range("K630").Value = 3800
$ws.Range("L630").Value = 4000
$ws.Range("M630").Value = 3900
$ws.Range("N630").Value = "$/paquete 36 unidades"
$ws.Range("O630").Value = "Provincia de Quillota"
$ws.Range("P630").Value = 108
$ws.Range("Q630").Value = 36
$ws.Range("R630").Value = "Hortaliza"
